# Applies the "Updated symbol list" edit: refreshed prices/hour-stamp
# for existing rows, plus a one-row insertion shift (ProBitToken at
# row 12 through BitpandaEcosystemToken at row 27) and a KickToken
# insertion shift (rows 41-43).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: BNB
$ws.Range('D2').Value = "'246.63"
$ws.Range('G2').Value = "'10"
# Row 3: OKB
$ws.Range('D3').Value = "'22.41"
$ws.Range('G3').Value = "'10"
# Row 4: HuobiToken
$ws.Range('D4').Value = "'5.500"
$ws.Range('G4').Value = "'10"
# Row 5: Cronos
$ws.Range('D5').Value = "'0.05614"
$ws.Range('G5').Value = "'10"
# Row 6: KuCoinToken
$ws.Range('D6').Value = "'6.468"
$ws.Range('G6').Value = "'10"
# Row 7: MXToken
$ws.Range('D7').Value = "'0.8017"
$ws.Range('G7').Value = "'10"
# Row 8: FTXToken
$ws.Range('D8').Value = "'1.046"
$ws.Range('G8').Value = "'10"
# Row 9: WazirX
$ws.Range('D9').Value = "'0.1423"
$ws.Range('G9').Value = "'10"
# Row 10: MandalaExchangeToken
$ws.Range('D10').Value = "'0.07248"
$ws.Range('G10').Value = "'10"
# Row 11: LiechtensteinCryptoassetsExchange
$ws.Range('D11').Value = "'0.03210"
$ws.Range('G11').Value = "'10"
# Row 12: ProBitToken
$ws.Range('B12').Value = 'ProBitToken'
$ws.Range('C12').Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range('D12').Value = "'0.1293"
$ws.Range('E12').Value = '11ProBitTokenPROBBestin24h'
$ws.Range('G12').Value = "'10"
# Row 13: BitrueCoin
$ws.Range('B13').Value = 'BitrueCoin'
$ws.Range('C13').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D13').Value = "'0.02969"
$ws.Range('E13').Value = '12BitrueCoinBTR'
$ws.Range('G13').Value = "'10"
# Row 14: BitMartToken
$ws.Range('B14').Value = 'BitMartToken'
$ws.Range('C14').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D14').Value = "'0.09262"
$ws.Range('E14').Value = '13BitMartTokenBMX'
$ws.Range('G14').Value = "'10"
# Row 15: BitForexToken
$ws.Range('B15').Value = 'BitForexToken'
$ws.Range('C15').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D15').Value = "'0.001663"
$ws.Range('E15').Value = '14BitForexTokenBF'
$ws.Range('G15').Value = "'10"
# Row 16: MCDex
$ws.Range('B16').Value = 'MCDex'
$ws.Range('C16').Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range('D16').Value = "'2.969"
$ws.Range('E16').Value = '15MCDexMCB'
$ws.Range('G16').Value = "'10"
# Row 17: CoinExToken
$ws.Range('B17').Value = 'CoinExToken'
$ws.Range('C17').Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range('D17').Value = "'0.04677"
$ws.Range('E17').Value = '16CoinExTokenCET'
$ws.Range('G17').Value = "'10"
# Row 18: One
$ws.Range('B18').Value = 'One'
$ws.Range('C18').Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range('D18').Value = "'0.0005983"
$ws.Range('E18').Value = '17OneONE'
$ws.Range('G18').Value = "'10"
# Row 19: TigerCash
$ws.Range('B19').Value = 'TigerCash'
$ws.Range('C19').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('D19').Value = "'0.006283"
$ws.Range('E19').Value = '18TigerCashTCH'
$ws.Range('G19').Value = "'10"
# Row 20: BitKan
$ws.Range('B20').Value = 'BitKan'
$ws.Range('C20').Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$ws.Range('D20').Value = "'0.001052"
$ws.Range('E20').Value = '19BitKanKAN'
$ws.Range('G20').Value = "'10"
# Row 21: HotbitToken
$ws.Range('B21').Value = 'HotbitToken'
$ws.Range('C21').Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Range('D21').Value = "'0.003810"
$ws.Range('E21').Value = '20HotbitTokenHTB'
$ws.Range('G21').Value = "'10"
# Row 22: NitroEx
$ws.Range('B22').Value = 'NitroEx'
$ws.Range('C22').Value = 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
$ws.Range('D22').Value = "'0.0001504"
$ws.Range('E22').Value = '21NitroExNTX'
$ws.Range('G22').Value = "'10"
# Row 23: UpBots
$ws.Range('B23').Value = 'UpBots'
$ws.Range('C23').Value = 'https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt'
$ws.Range('D23').Value = "'0.0003608"
$ws.Range('E23').Value = '22UpBotsUBXT'
$ws.Range('G23').Value = "'10"
# Row 24: LEO
$ws.Range('B24').Value = 'LEO'
$ws.Range('C24').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D24').Value = "'3.985"
$ws.Range('E24').Value = '23LEOLEO'
$ws.Range('G24').Value = "'10"
# Row 25: GateToken
$ws.Range('B25').Value = 'GateToken'
$ws.Range('C25').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range('D25').Value = "'3.400"
$ws.Range('E25').Value = '24GateTokenGT'
$ws.Range('G25').Value = "'10"
# Row 26: BTSEToken
$ws.Range('B26').Value = 'BTSEToken'
$ws.Range('C26').Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range('D26').Value = "'2.122"
$ws.Range('E26').Value = '25BTSETokenBTSE'
$ws.Range('G26').Value = "'10"
# Row 27: BitpandaEcosystemToken
$ws.Range('B27').Value = 'BitpandaEcosystemToken'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range('D27').Value = "'0.3290"
$ws.Range('E27').Value = '26BitpandaEcosystemTokenBEST'
$ws.Range('G27').Value = "'10"
# Row 28: Spectre.aiUtilityToken
$ws.Range('G28').Value = "'10"
# Row 29: LegolasExchange
$ws.Range('G29').Value = "'10"
# Row 30: BitZToken
$ws.Range('G30').Value = "'10"
# Row 31: Birake
$ws.Range('G31').Value = "'10"
# Row 32: ZBToken
$ws.Range('G32').Value = "'10"
# Row 33: NashExchange
$ws.Range('G33').Value = "'10"
# Row 34: AAXToken
$ws.Range('G34').Value = "'10"
# Row 35: CenX
$ws.Range('G35').Value = "'10"
# Row 36: BNIXToken
$ws.Range('G36').Value = "'10"
# Row 37: Polkally
$ws.Range('G37').Value = "'10"
# Row 38: Charli3
$ws.Range('G38').Value = "'10"
# Row 39: BlubitexToken
$ws.Range('G39').Value = "'10"
# Row 40: IDEX
$ws.Range('D40').Value = "'0.04146"
$ws.Range('G40').Value = "'10"
# Row 41: KickToken
$ws.Range('B41').Value = 'KickToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range('D41').Value = "'0.007006"
$ws.Range('E41').Value = '40KickTokenKICK'
$ws.Range('G41').Value = "'10"
# Row 42: BKEXToken
$ws.Range('B42').Value = 'BKEXToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range('D42').Value = "'0.1038"
$ws.Range('E42').Value = '41BKEXTokenBKK'
$ws.Range('G42').Value = "'10"
# Row 43: CEJI
$ws.Range('B43').Value = 'CEJI'
$ws.Range('C43').Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range('D43').Value = "'0.003163"
$ws.Range('E43').Value = '42CEJICEJI'
$ws.Range('G43').Value = "'10"
# Row 44: LocalTraders
$ws.Range('D44').Value = "'0.01027"
$ws.Range('G44').Value = "'10"
# Row 45: CoinLion
$ws.Range('D45').Value = "'0.00005655"
$ws.Range('G45').Value = "'10"
# Row 46: Kangarootoken
$ws.Range('G46').Value = "'10"
# Row 47: CoinbaseStockToken
$ws.Range('D47').Value = "'0.6816"
$ws.Range('G47').Value = "'10"
# Row 48: BOLO
$ws.Range('D48').Value = "'0.02636"
$ws.Range('E48').Value = '47BOLOBOLOWorstin24h'
$ws.Range('G48').Value = "'10"
# Row 49: CryptobidCoin
$ws.Range('D49').Value = "'0.00002105"
$ws.Range('G49').Value = "'10"
# Row 50: SpecialPowerGold
$ws.Range('D50').Value = "'0.01012"
$ws.Range('G50').Value = "'10"
# Row 51: DigiFinexToken
$ws.Range('G51').Value = "'10"
